$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data with the "Graduate" / "GRAD" class level
$ws.Range("A11").Value = "Graduate"
$ws.Range("C11").Value = "GRAD"
$ws.Range("D11").Value = 21
$ws.Range("D11").NumberFormat = "`"$`"#,##0_);[Red]\(`"$`"#,##0\)"
$ws.Range("E11").Value = "60 min"
$ws.Range("H11").Value = 6
$ws.Range("I11").Value = 12

$ws.Range("C12").Select()
